$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update value cells (non-formula)
$ws.Range("B3").Value = 0
$ws.Range("B10").Value = 1.2
$ws.Range("B15").Value = 1.5

# Formulas in B7 and B8 depend on B3, so they will recalculate automatically.

# Update the selection to match the saved view state (active cell F6)
$ws.Range("F6").Select()
